# Refresh cryptocurrency price (D) and 1h volume-change (E) figures on Sheet1,
# matching the latest scrape from the GitHub Actions update job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.969.06"
$ws.Range("E2").Value = "  +2.82%  "
$ws.Range("D3").Value = "3.806.69"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'702.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +11.55%  "
$ws.Range("D6").Value = "'173.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.77%  "
$ws.Range("D7").Value = "3.806.15"
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("D10").Value = "'0.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.21%  "
$ws.Range("D11").Value = "'7.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.54%  "
$ws.Range("E12").Value = "  +1.68%  "
$ws.Range("D13").Value = "'0.0000261"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.68%  "
$ws.Range("D14").Value = "'36.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.43%  "
$ws.Range("D15").Value = "4.446.51"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D16").Value = "3.806.55"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "70.997.39"
$ws.Range("E17").Value = "  +2.84%  "
$ws.Range("D18").Value = "'17.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("D19").Value = "'7.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.26%  "
$ws.Range("D21").Value = "'11.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +17.50%  "
$ws.Range("D22").Value = "'484.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.36%  "
$ws.Range("D23").Value = "'0.716"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("D24").Value = "'84.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.59%  "
$ws.Range("E25").Value = "  +3.35%  "
$ws.Range("D26").Value = "'12.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'10.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.07%  "
$ws.Range("D29").Value = "3.957.78"
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'3.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +15.66%  "
$ws.Range("D32").Value = "'7.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.96%  "
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("D34").Value = "'29.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.39%  "
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("D36").Value = "'9.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.44%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "3.757.37"
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("D39").Value = "'0.103"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.70%  "
$ws.Range("D40").Value = "'3.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.44%  "
$ws.Range("D41").Value = "'5.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.61%  "
$ws.Range("D42").Value = "'2.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.45%  "
$ws.Range("D43").Value = "'0.000329"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +24.93%  "
$ws.Range("D44").Value = "'0.972"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "'162.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.06%  "
$ws.Range("D48").Value = "'49.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.28%  "
$ws.Range("D49").Value = "'44.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.43%  "
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").Value = "'1.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.79%  "
